{"js": "// Applies the \"With the current updates\" revision to The Music Notes.docx.\n//\n// Two kinds of changes:\n//   (A) A handful of existing paragraphs had their Word-proofing marks\n//       (w:proofErr spellStart/spellEnd/gramStart/gramEnd) cleared away and\n//       the runs that were only split apart to bracket those marks were\n//       rejoined. The visible text of every one of these paragraphs is\n//       unchanged - only the run/proofErr bookkeeping differs.\n//   (B) A large block of new meeting-notes paragraphs was appended at the\n//       end of the document (after \"Data by year, month, etc.\").\n//\n// We use Range.insertOoxml(...) throughout so the resulting run/proofErr\n// structure matches the target exactly instead of relying on how the\n// higher-level text APIs happen to merge runs.\n\nconst OOXML_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapPackage(bodyInnerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    `<w:document ${OOXML_NS}>` +\n    \"<w:body>\" +\n    bodyInnerXml +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nfunction replaceParagraphWith(paragraph, innerBodyXml) {\n  const range = paragraph.getRange(Word.RangeLocation.whole);\n  range.insertOoxml(wrapPackage(innerBodyXml), Word.InsertLocation.replace);\n}\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\nconst paras = body.paragraphs.items;\n\n// ---------------------------------------------------------------------\n// (A) Proofing-mark cleanup / run merges on existing paragraphs.\n// ---------------------------------------------------------------------\n\n// \"Team Members: Michelle Petras, Veethika Singh, Kathleen Snider-Belinski, Carmen Wiggins\"\nreplaceParagraphWith(\n  paras[2],\n  \"<w:p><w:r><w:t>Team Members: Michelle Petras, Veethika Singh, Kathleen Snider-Belinski, Carmen Wiggins</w:t></w:r></w:p>\"\n);\n\n// \"Datasets to be used:  Kaggle Spotify datasets, and possible API calls\"\nreplaceParagraphWith(\n  paras[6],\n  \"<w:p><w:r><w:t>Datasets to be used:  Kaggle Spotify datasets, and possible API calls</w:t></w:r></w:p>\"\n);\n\n// \"Create initial Git setup \u2013 Michelle\"\nreplaceParagraphWith(\n  paras[10],\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    \"<w:r><w:t>Create initial Git setup \\u2013 Michelle</w:t></w:r></w:p>\"\n);\n\n// \"Clone Git setup \u2013 Veethika, Kat, Carmen\"\nreplaceParagraphWith(\n  paras[11],\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    \"<w:r><w:t>Clone Git setup \\u2013 Veethika, Kat, Carmen</w:t></w:r></w:p>\"\n);\n\n// \"Download data set(s) - ?\"\nreplaceParagraphWith(\n  paras[12],\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n    \"<w:r><w:t>Download data set(s) - ?</w:t></w:r></w:p>\"\n);\n\n// \"precovid\" + \" 01/01/2017 \u2013 03/17/2020\" -- proofErr removed, runs stay split.\nreplaceParagraphWith(\n  paras[16],\n  \"<w:p><w:r><w:t>precovid</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\"> 01/01/2017 \\u2013 03/17/2020</w:t></w:r></w:p>'\n);\n\n// \"postcovid 03/18/2020 - 2022\" -- proofErr removed AND runs merged into one.\nreplaceParagraphWith(paras[17], \"<w:p><w:r><w:t>postcovid 03/18/2020 - 2022</w:t></w:r></w:p>\");\n\n// \"Get edit csv\"\nreplaceParagraphWith(paras[19], \"<w:p><w:r><w:t>Get edit csv</w:t></w:r></w:p>\");\n\n// \"March 17\" + superscript \"th\" stay untouched; the trailing two runs merge\n// into \" 2020 is demarcation point.\" with the gramStart/gramEnd marks gone.\nreplaceParagraphWith(\n  paras[22],\n  \"<w:p><w:r><w:t>March 17</w:t></w:r>\" +\n    '<w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>th</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> 2020 is demarcation point.</w:t></w:r></w:p>'\n);\n\n// ---------------------------------------------------------------------\n// (B) New paragraphs appended after \"Data by year, month, etc.\"\n// ---------------------------------------------------------------------\n\nconst newParagraphsXml =\n  \"<w:p/>\" +\n  \"<w:p/>\" +\n  \"<w:p><w:r><w:t>Cleanup genre - Kat</w:t></w:r></w:p>\" +\n  \"<w:p>\" +\n  \"<w:r><w:t>Streams by genre</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n  \"<w:r><w:t>yearly</w:t></w:r>\" +\n  \"<w:r><w:t>?  Monthly?</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  \"<w:r><w:t>Veethika</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr><w:strike/></w:rPr><w:t>NY times API articles?</w:t></w:r>\" +\n  '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space=\"preserve\"> - </w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:r><w:t>Total number of streams</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> - </w:t></w:r>' +\n  \"<w:r><w:t>Carmen</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:r><w:lastRenderedPageBreak/><w:t>Top performing artists, pre and post</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> - Carmen</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:r><w:t>Loudness</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  \"<w:r><w:t>(Db)</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> Michelle</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:r><w:t>Tempo (bpm)</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> Michelle</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:r><w:t>Danceability</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> Michelle</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:r><w:t>Speechiness</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> Michelle</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:r><w:t>Energy</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  \"<w:r><w:t>(0-1)</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> Michelle</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr><w:strike/></w:rPr><w:t>By region or country?</w:t></w:r>\" +\n  '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space=\"preserve\"> - Veethika</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p>\" +\n  \"<w:r><w:rPr><w:strike/></w:rPr><w:t># of users?</w:t></w:r>\" +\n  '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space=\"preserve\"> - Veethika</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p><w:r><w:t>Is it # of times it is streamed or by unique users?</w:t></w:r></w:p>\";\n\nconst lastParagraph = paras[paras.length - 1];\nconst endRange = lastParagraph.getRange(Word.RangeLocation.end);\nendRange.insertOoxml(wrapPackage(newParagraphsXml), Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Applies the \"With the current updates\" revision to The Music Notes.docx.\n#\n# Two kinds of changes:\n#   (A) A handful of existing paragraphs had their Word-proofing marks\n#       (w:proofErr spellStart/spellEnd/gramStart/gramEnd) cleared away and\n#       the runs that were only split apart to bracket those marks were\n#       rejoined. The visible text of every one of these paragraphs is\n#       unchanged - only the run/proofErr bookkeeping differs.\n#   (B) A large block of new meeting-notes paragraphs was appended at the\n#       end of the document (after \"Data by year, month, etc.\").\n#\n# We drive this through Range.InsertXML(...) (WordprocessingML package XML)\n# so the resulting run/proofErr structure matches the target exactly,\n# instead of relying on how Range.Text assignment happens to merge runs.\n\n$d = $word.ActiveDocument\n\nfunction Wrap-Package($bodyInnerXml) {\n    $head = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>'\n    $tail = '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n    $full = $head + $bodyInnerXml + $tail\n    return $full\n}\n\nfunction Set-ParagraphOoxml($paraIndex1Based, $innerBodyXml) {\n    $p = $d.Paragraphs($paraIndex1Based)\n    $r = $p.Range\n    $pkg = Wrap-Package $innerBodyXml\n    $r.InsertXML($pkg)\n}\n\n$enDash = [char]0x2013\n\n# ---------------------------------------------------------------------\n# (A) Proofing-mark cleanup / run merges on existing paragraphs.\n# (Paragraphs() is 1-indexed, so body-child index N -> Paragraphs(N+1).)\n# ---------------------------------------------------------------------\n\n# Paragraph 3: \"Team Members: Michelle Petras, Veethika Singh, Kathleen Snider-Belinski, Carmen Wiggins\"\n$p3 = '<w:p><w:r><w:t>Team Members: Michelle Petras, Veethika Singh, Kathleen Snider-Belinski, Carmen Wiggins</w:t></w:r></w:p>'\nSet-ParagraphOoxml 3 $p3\n\n# Paragraph 7: \"Datasets to be used:  Kaggle Spotify datasets, and possible API calls\"\n$p7 = '<w:p><w:r><w:t>Datasets to be used:  Kaggle Spotify datasets, and possible API calls</w:t></w:r></w:p>'\nSet-ParagraphOoxml 7 $p7\n\n$listPPr = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>'\n\n# Paragraph 11: \"Create initial Git setup - Michelle\" (en dash)\n$p11 = '<w:p>' + $listPPr + '<w:r><w:t>Create initial Git setup ' + $enDash + ' Michelle</w:t></w:r></w:p>'\nSet-ParagraphOoxml 11 $p11\n\n# Paragraph 12: \"Clone Git setup - Veethika, Kat, Carmen\" (en dash)\n$p12 = '<w:p>' + $listPPr + '<w:r><w:t>Clone Git setup ' + $enDash + ' Veethika, Kat, Carmen</w:t></w:r></w:p>'\nSet-ParagraphOoxml 12 $p12\n\n# Paragraph 13: \"Download data set(s) - ?\"\n$p13 = '<w:p>' + $listPPr + '<w:r><w:t>Download data set(s) - ?</w:t></w:r></w:p>'\nSet-ParagraphOoxml 13 $p13\n\n# Paragraph 17: \"precovid\" stays split into two runs; just the proofErr goes away.\n$p17 = '<w:p><w:r><w:t>precovid</w:t></w:r><w:r><w:t xml:space=\"preserve\"> 01/01/2017 ' + $enDash + ' 03/17/2020</w:t></w:r></w:p>'\nSet-ParagraphOoxml 17 $p17\n\n# Paragraph 18: \"postcovid 03/18/2020 - 2022\" merges into a single run.\n$p18 = '<w:p><w:r><w:t>postcovid 03/18/2020 - 2022</w:t></w:r></w:p>'\nSet-ParagraphOoxml 18 $p18\n\n# Paragraph 20: \"Get edit csv\" merges into a single run.\n$p20 = '<w:p><w:r><w:t>Get edit csv</w:t></w:r></w:p>'\nSet-ParagraphOoxml 20 $p20\n\n# Paragraph 23: \"March 17th 2020 is demarcation point.\" -- keep \"March 17\" and\n# the superscript \"th\" run untouched; merge the trailing two runs.\n$p23 = '<w:p><w:r><w:t>March 17</w:t></w:r>' +\n    '<w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>th</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> 2020 is demarcation point.</w:t></w:r></w:p>'\nSet-ParagraphOoxml 23 $p23\n\n# ---------------------------------------------------------------------\n# (B) New paragraphs appended after \"Data by year, month, etc.\"\n# ---------------------------------------------------------------------\n\n$strikeRPr = '<w:rPr><w:strike/></w:rPr>'\n\n$paraCleanupGenre = '<w:p><w:r><w:t>Cleanup genre - Kat</w:t></w:r></w:p>'\n\n$paraStreamsByGenre = '<w:p>' +\n        '<w:r><w:t>Streams by genre</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n        '<w:r><w:t>yearly</w:t></w:r>' +\n        '<w:r><w:t>?  Monthly?</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n        '<w:r><w:t>Veethika</w:t></w:r>' +\n    '</w:p>'\n\n$paraNyTimes = '<w:p>' +\n        '<w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>' +\n        '<w:r>' + $strikeRPr + '<w:t>NY times API articles?</w:t></w:r>' +\n        '<w:r>' + $strikeRPr + '<w:t xml:space=\"preserve\"> - </w:t></w:r>' +\n    '</w:p>'\n\n$paraTotalStreams = '<w:p>' +\n        '<w:r><w:t>Total number of streams</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> - </w:t></w:r>' +\n        '<w:r><w:t>Carmen</w:t></w:r>' +\n    '</w:p>'\n\n$paraTopArtists = '<w:p>' +\n        '<w:r><w:lastRenderedPageBreak/><w:t>Top performing artists, pre and post</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> - Carmen</w:t></w:r>' +\n    '</w:p>'\n\n$paraLoudness = '<w:p>' +\n        '<w:r><w:t>Loudness</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n        '<w:r><w:t>(Db)</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> Michelle</w:t></w:r>' +\n    '</w:p>'\n\n$paraTempo = '<w:p>' +\n        '<w:r><w:t>Tempo (bpm)</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> Michelle</w:t></w:r>' +\n    '</w:p>'\n\n$paraDanceability = '<w:p>' +\n        '<w:r><w:t>Danceability</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> Michelle</w:t></w:r>' +\n    '</w:p>'\n\n$paraSpeechiness = '<w:p>' +\n        '<w:r><w:t>Speechiness</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> Michelle</w:t></w:r>' +\n    '</w:p>'\n\n$paraEnergy = '<w:p>' +\n        '<w:r><w:t>Energy</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n        '<w:r><w:t>(0-1)</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> Michelle</w:t></w:r>' +\n    '</w:p>'\n\n$paraRegion = '<w:p>' +\n        '<w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>' +\n        '<w:r>' + $strikeRPr + '<w:t>By region or country?</w:t></w:r>' +\n        '<w:r>' + $strikeRPr + '<w:t xml:space=\"preserve\"> - Veethika</w:t></w:r>' +\n    '</w:p>'\n\n$paraUsers = '<w:p>' +\n        '<w:r>' + $strikeRPr + '<w:t># of users?</w:t></w:r>' +\n        '<w:r>' + $strikeRPr + '<w:t xml:space=\"preserve\"> - Veethika</w:t></w:r>' +\n    '</w:p>'\n\n$paraUniqueUsers = '<w:p><w:r><w:t>Is it # of times it is streamed or by unique users?</w:t></w:r></w:p>'\n\n$newParagraphsXml = '<w:p/>' + '<w:p/>' +\n    $paraCleanupGenre +\n    $paraStreamsByGenre +\n    $paraNyTimes +\n    $paraTotalStreams +\n    $paraTopArtists +\n    $paraLoudness +\n    $paraTempo +\n    $paraDanceability +\n    $paraSpeechiness +\n    $paraEnergy +\n    $paraRegion +\n    $paraUsers +\n    $paraUniqueUsers\n\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endPkg = Wrap-Package $newParagraphsXml\n$endRange.InsertXML($endPkg)\n\n\"done\"\n"}
